$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Cutoff analysis" entry (old row 5) moves down to become the new
# row 7, the "Covid visualisation..." entry (old row 6) shifts up to
# row 5, and a brand-new row 6 is inserted for the new blog post.

# New row 5 <- old row 6 (Covid visualisation using leaflet)
$ws.Range("B5").Value = "Covid visualisation using leaflet"
$ws.Range("C5").Value = "Data visualisation series"
$ws.Range("D5").Value = 2
$ws.Range("E5").ClearContents()

# New row 6 <- brand new post about XGBoost monotonic binning
$ws.Range("B6").Value = "Monotonic binning using XGBOOST"
$ws.Range("C6").Value = "Credit risk series"
$ws.Range("D6").Value = 2

# New row 7 <- old row 5 (Cutoff analysis), now without the D value
$ws.Range("B7").Value = "Cutoff analysis"
$ws.Range("C7").Value = "Credit risk series"
$ws.Range("D7").ClearContents()
$ws.Range("E7").Value = "Incomplete"

# Update the active selection to match the saved workbook state.
$ws.Range("C6").Select()
